# Update cryptocurrency price/volume data in Sheet1 (columns D and E, rows 2-51)
# to reflect the latest scrape, per the GitHub Actions automated commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.335.21"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "1.862.66"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "234.09"
$ws.Range("E5").Value = "  -2.24%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "0.4758"
$ws.Range("E7").Value = "  -1.12%  "
$ws.Range("D8").Value = "0.2756"
$ws.Range("E8").Value = "  -2.93%  "
$ws.Range("D9").Value = "0.06448"
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D10").Value = "1.850.83"
$ws.Range("E10").Value = "  -13.67%  "
$ws.Range("D11").Value = "0.07428"
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("D12").Value = "16.15"
$ws.Range("E12").Value = "  -3.20%  "
$ws.Range("D13").Value = "4.996"
$ws.Range("E13").Value = "  -2.19%  "
$ws.Range("D14").Value = "86.04"
$ws.Range("E14").Value = "  -3.15%  "
$ws.Range("D15").Value = "0.6342"
$ws.Range("E15").Value = "  -4.77%  "
$ws.Range("D16").Value = "30.300.49"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Value = "0.9996"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "232.47"
$ws.Range("E18").Value = "  +3.21%  "
$ws.Range("E19").Value = "  -4.08%  "
$ws.Range("D20").Value = "0.000007391"
$ws.Range("E20").Value = "  -3.03%  "
$ws.Range("D21").Value = "2.096.97"
$ws.Range("E21").Value = "  -4.11%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "5.105"
$ws.Range("E23").Value = "  -4.40%  "
$ws.Range("D24").Value = "6.029"
$ws.Range("E24").Value = "  -3.31%  "
$ws.Range("D25").Value = "9.307"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("D26").Value = "167.37"
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("D27").Value = "17.94"
$ws.Range("E27").Value = "  -3.86%  "
$ws.Range("D28").Value = "1.864"
$ws.Range("E28").Value = "  -5.07%  "
$ws.Range("D29").Value = "1.381"
$ws.Range("E29").Value = "  -5.41%  "
$ws.Range("D30").Value = "0.1004"
$ws.Range("E30").Value = "  +5.70%  "
$ws.Range("D31").Value = "4.233"
$ws.Range("E31").Value = "  -2.38%  "
$ws.Range("D32").Value = "3.919"
$ws.Range("E32").Value = "  -3.23%  "
$ws.Range("E33").Value = "  -2.68%  "
$ws.Range("D34").Value = "1.152"
$ws.Range("E34").Value = "  -4.85%  "
$ws.Range("D35").Value = "0.7257"
$ws.Range("E35").Value = "  -3.66%  "
$ws.Range("D36").Value = "0.9990"
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("D37").Value = "2.692"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("D38").Value = "0.01937"
$ws.Range("E38").Value = "  +5.40%  "
$ws.Range("D39").Value = "2.633"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").Value = "0.9093"
$ws.Range("D41").Value = "1.991"
$ws.Range("E41").Value = "  -4.35%  "
$ws.Range("D42").Value = "105.63"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").Value = "0.4117"
$ws.Range("E44").Value = "  -4.24%  "
$ws.Range("D45").Value = "5.547"
$ws.Range("E45").Value = "  -4.73%  "
$ws.Range("D46").Value = "7.082"
$ws.Range("E46").Value = "  -5.62%  "
$ws.Range("D47").Value = "61.39"
$ws.Range("E47").Value = "  -5.56%  "
$ws.Range("E48").Value = "  -6.28%  "
$ws.Range("D49").Value = "8.770"
$ws.Range("E49").Value = "  -2.56%  "
$ws.Range("D50").Value = "1.403"
$ws.Range("E50").Value = "  -4.99%  "
$ws.Range("D51").Value = "33.11"
$ws.Range("E51").Value = "  -2.42%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"

Write-Host "Updated cryptos list"
